# Workbook/worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the duplicated header row (old row 1, "Unnamed: 0"/bold
#    bordered style). This shifts every row up by one, so the old
#    row 2 (plain header) becomes the new row 1, old row 3
#    ("Fixation based metrics") becomes row 2, etc. It also drops the
#    bold/border style that was only applied to that first row.
# ------------------------------------------------------------------
$ws.Rows("1:1").Delete()

# ------------------------------------------------------------------
# 2. Correct the data-cleaning values for the metric rows (now rows
#    3-8 after the shift above).
# ------------------------------------------------------------------

# Row 3: Revisit count
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 2
$ws.Range("F3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 15
$ws.Range("O3").Value = 0

# Row 4: Fixation count
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 59
$ws.Range("O4").Value = 1

# Row 5: Dwell time (ms)
$ws.Range("B5").Value = 216.02
$ws.Range("C5").Value = 849.89
$ws.Range("D5").Value = 367.01
$ws.Range("F5").Value = 150.16
$ws.Range("K5").Value = 216.02
$ws.Range("L5").Value = 583.03
$ws.Range("M5").Value = 13747.97
$ws.Range("O5").Value = 216.02

# Row 6: Dwell time (%)
$ws.Range("B6").Value = 0.18
$ws.Range("C6").Value = 0.6899999999999999
$ws.Range("D6").Value = 0.3
$ws.Range("F6").Value = 0.12
$ws.Range("K6").Value = 0.18
$ws.Range("L6").Value = 0.47
$ws.Range("M6").Value = 11.18
$ws.Range("O6").Value = 0.18

# Row 7: Fixation duration (ms)
$ws.Range("B7").Value = 216.02
$ws.Range("C7").Value = 169.98
$ws.Range("D7").Value = 183.51
$ws.Range("F7").Value = 150.16
$ws.Range("K7").Value = 216.02
$ws.Range("L7").Value = 194.34
$ws.Range("M7").Value = 233.02
$ws.Range("O7").Value = 216.02

# Row 8: First fixation duration (ms)
$ws.Range("B8").Value = 216.02
$ws.Range("C8").Value = 216.02
$ws.Range("D8").Value = 150.16
$ws.Range("F8").Value = 150.16
$ws.Range("K8").Value = 216.02
$ws.Range("L8").Value = 216.02
$ws.Range("M8").Value = 150.11
$ws.Range("O8").Value = 216.02
